$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.078222337956418642
$ws.Range("A2").Value = -0.0099999997733988266
$ws.Range("A3").Value = -0.0089999997769059092
$ws.Range("A4").Value = 0.28399534826292694
$ws.Range("A5").Value = -0.0059999997861828192
$ws.Range("A6").Value = -0.0059999997807480554
$ws.Range("A7").Value = -0.019999999738869789
$ws.Range("A8").Value = -0.019999999736896257
$ws.Range("A9").Value = -0.0059999997762254509
$ws.Range("A10").Value = -0.005999999774722653
$ws.Range("A11").Value = -0.0044999997794796798
$ws.Range("A12").Value = 0.0058937587193352847
$ws.Range("A13").Value = -0.0059999997731692289
$ws.Range("A14").Value = -0.011999999755131441
$ws.Range("A15").Value = -0.0059999997726984944
$ws.Range("A16").Value = 0.019113679549294105
$ws.Range("A17").Value = -0.0059999997718032105
$ws.Range("A18").Value = -0.0089999997627048245
$ws.Range("A19").Value = -0.0089999997755532135
$ws.Range("A20").Value = -0.070676826763985012
$ws.Range("A21").Value = -0.0089999997678029686
$ws.Range("A22").Value = -0.008999999767531186
$ws.Range("A23").Value = -0.0089999997742031823
$ws.Range("A24").Value = -0.06695028898810218
$ws.Range("A25").Value = -0.041999999670036203
$ws.Range("A26").Value = -0.0059999997804709437
$ws.Range("A27").Value = -0.0059999997800526117
$ws.Range("A28").Value = -0.0059999997785666892
$ws.Range("A29").Value = -0.0119999997600857
$ws.Range("A30").Value = 0.033820654350907642
$ws.Range("A31").Value = -0.014999999748877002
$ws.Range("A32").Value = -0.020999999730807239
$ws.Range("A33").Value = -0.0059999997748310108
